$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1D NEW")

# --- Add three new simulated-interferogram rows to the bottom of the table ---

# Row 26: 1dmockanderrors23.csv
$ws.Range("B26").Value = "1dmockanderrors23.csv"
$ws.Range("C26").Value = 53
$ws.Range("D26").Value = 1000
$ws.Range("E26").Value = 0.27
$ws.Range("F26").Value = 0.05
$ws.Range("G26").Value = 200
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 5
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = "we're so back"

# Row 27: 1dmockanderrors24.csv (comment entered before the filename)
$ws.Range("C27").Value = 53
$ws.Range("D27").Value = 1000
$ws.Range("E27").Value = 0.27
$ws.Range("F27").Value = 0.05
$ws.Range("G27").Value = 200
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 10
$ws.Range("L27").Value = 1
$ws.Range("M27").Value = "same but with more noise"
$ws.Range("B27").Value = "1dmockanderrors24.csv"

# Row 28: 1dmockanderrors25.csv
$ws.Range("B28").Value = "1dmockanderrors25.csv"
$ws.Range("C28").Value = 53
$ws.Range("D28").Value = 1000
$ws.Range("E28").Value = 0.27
$ws.Range("F28").Value = 0.1
$ws.Range("G28").Value = 200
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 5
$ws.Range("L28").Value = 1
$ws.Range("M28").Value = "same but with larger FWHM"

# --- Grow the "Table43" structured table to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B3:M28")) | Out-Null

# --- Re-made graph off new data means this sheet is now the one in focus ---
$ws.Activate() | Out-Null
$ws.Range("F29").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
